$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Tacos" restaurant entries appended below the existing data (rows 42-49)
$newRows = @(
    @("Taco Joint",        "Tacos", "$$",   "Central"),
    @("El Chilito",        "Tacos", "$$",   "East"),
    @("Juan in a Million",  "Tacos", "$",    "East"),
    @("Taco Flats",        "Tacos", "$$$",  "Central"),
    @("Vaquero Taquero",   "Tacos", "$$",   "Central"),
    @("Veracruz",          "Tacos", "$$",   "Central"),
    @("Rudy's",            "Tacos", "$$",   "Central"),
    @("Texas Honey Ham",   "Tacos", "$$",   "West")
)

$startRow = 42

# Shared strings are interned in the order cells are written, matching the
# source data's column-major entry order (all of column A, then B, then C,
# then D) so new unique strings land in the same sequence as the target file.
for ($col = 1; $col -le 4; $col++) {
    for ($i = 0; $i -lt $newRows.Count; $i++) {
        $row = $startRow + $i
        $data = $newRows[$i]
        $ws.Cells.Item($row, $col).Value = $data[$col - 1]
    }
}

$ws.Range("D50").Select()
